$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Correct the release id text (A3): "COMMANDCARD V1.0" -> "COMCRD  v1.0"
$ws.Range("A3").Value = "COMCRD  v1.0"

# 2. Move active selection from I3 to A3
$ws.Range("A3").Select() | Out-Null

# 3. Widen columns A:J (columns were resized wider by roughly the same proportion)
$ws.Range("A:A").ColumnWidth = 22.358258928571427
$ws.Range("B:B").ColumnWidth = 44.07254464285714
$ws.Range("C:C").ColumnWidth = 21.50111607142857
$ws.Range("D:H").ColumnWidth = 13.072544642857142
$ws.Range("I:I").ColumnWidth = 32.9296875
$ws.Range("J:J").ColumnWidth = 11.9296875

# 4. Shrink the print scale from 100% to 55%
$ws.PageSetup.Zoom = 55

# 5. Update the header text from LEDCON.MECH.ECO to COMCRD.ECO
$ws.PageSetup.LeftHeader = "&""Times New Roman,Regular""COMCRD.ECO"
